$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet_cut_off")
$ws.Activate()

# Widen column A to fit the new content
$ws.Columns.Item(1).ColumnWidth = 16

# Fill in the numeric cut-off values for rows 2-5 (columns B-E)
$ws.Range("B2").Value = 23
$ws.Range("C2").Value = 32
$ws.Range("D2").Value = 39
$ws.Range("E2").Value = 40

$ws.Range("B3").Value = 12
$ws.Range("C3").Value = 17
$ws.Range("D3").Value = 20
$ws.Range("E3").Value = 21

$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 6
$ws.Range("D4").Value = 10
$ws.Range("E4").Value = 11

$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 8
$ws.Range("E5").Value = 9

# Move the active selection to G6 as in the edited workbook
$ws.Range("G6").Select()
